# actualizacion 7 de Marzo
# Adds two new rows (20-21) of data to the "Resumen" sheet:
#   Row 20: "Número de transacciones" with three numeric totals
#   Row 21: "Valor transacciones"    with three large comma-grouped amounts
#           stored as TEXT (not numbers)
#
# The comma-grouped amounts (e.g. "6,607,677,129,791") must land in the
# sheet as literal text, exactly like the source workbook. A plain
# Range.Value/.Formula assignment of such a string gets auto-parsed into a
# number by Excel's "smart" literal-entry heuristics. To avoid that (and
# avoid introducing an extra, unwanted cell style such as quotePrefix or a
# "@" text format), we build the text with a TEXT() formula in a scratch
# cell, copy it, and paste-special *values only* into the destination -
# clipboard paste bypasses the literal-entry parser, landing a plain shared
# string with the default cell style, then we clean the scratch cell back up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: "Número de transacciones" -------------------------------
$ws.Range("A20").Value = "Número de transacciones"
$ws.Range("B20").Value = 19411918
$ws.Range("C20").Value = 22301752
$ws.Range("D20").Value = 41713670

# --- Row 21: "Valor transacciones" (text amounts) ---------------------
$ws.Range("A21").Value = "Valor transacciones"

$scratch = $ws.Range("Z1")

$scratch.Formula = "=TEXT(6607677129791,""#,##0"")"
$scratch.Copy()
$ws.Range("B21").PasteSpecial(-4163)

$scratch.Formula = "=TEXT(7421985990886,""#,##0"")"
$scratch.Copy()
$ws.Range("C21").PasteSpecial(-4163)

$scratch.Formula = "=TEXT(14029663120677,""#,##0"")"
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)

$scratch.ClearContents()
$excel.CutCopyMode = $false

# --- Match the author's final selection (A20:D21, active cell A20) ----
$ws.Range("A20:D21").Select()
